# "Store Plan for zero-doc"
#
# The ambient directory seed data renamed two store-plan JSON payload
# references (cab/directory/service.catalog.json and
# cab/directory/system.document.json) to the new zero-doc naming scheme
# (cab/directory/ambient.dir.document.json and
# cab/directory/ambient.dir.workflow.json) on the DATA-CDATA sheet, row 5
# (K5) and row 6 (K6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K5 held the "system.document.json" store plan -> now the document plan
# under the new ambient.dir.* naming.
$ws.Range("K5").Value = "JSON:cab/directory/ambient.dir.document.json"

# K6 held the "service.catalog.json" store plan -> now the workflow plan
# under the new ambient.dir.* naming.
$ws.Range("K6").Value = "JSON:cab/directory/ambient.dir.workflow.json"

# The new, longer text no longer fits the old best-fit width of column K,
# so it is widened to match.
$ws.Columns("K").ColumnWidth = 55.333333333333336

# Leave the active selection on K7, as in the saved workbook.
$ws.Range("K7").Select()
